# Update forecast values in the "Forecast Comparison" sheet (Removed Auto Arima).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$data = @{
    2  = @(70, 85, 106, 141)
    3  = @(68, 84, 103, 134)
    4  = @(98, 119, 142, 179)
    5  = @(98, 119, 143, 180)
    6  = @(109, 133, 162, 208)
    7  = @(106, 129, 156, 199)
    8  = @(108, 132, 162, 211)
    9  = @(109, 133, 163, 211)
    10 = @(104, 127, 154, 198)
    11 = @(107, 131, 161, 208)
    12 = @(110, 135, 166, 217)
    13 = @(113, 139, 172, 225)
    14 = @(112, 137, 169, 221)
    15 = @(113, 138, 172, 227)
    16 = @(111, 136, 169, 222)
    17 = @(105, 128, 161, 214)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
}
